$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add column D "Final self-report" values (BRIEF / CLASS responses)
$ws.Range("D2").Value = "Somewhat worse"
$ws.Range("D3").Value = "A lot worse"
$ws.Range("D4").Value = "A lot worse"
$ws.Range("D5").Value = "A lot worse"
$ws.Range("D6").Value = "A lot worse"
$ws.Range("D7").Value = "A lot worse"
$ws.Range("D8").Value = "A lot worse"
$ws.Range("D9").Value = "Somewhat worse"
$ws.Range("D10").Value = "Somewhat worse"
$ws.Range("D11").Value = "A lot worse"
$ws.Range("D12").Value = "A lot worse"
$ws.Range("D13").Value = "Somewhat worse"
$ws.Range("D14").Value = "A little worse"
$ws.Range("D15").Value = "A lot worse"

# Widen column D to match column B's width
$ws.Columns("D").ColumnWidth = 18.6640625

# Update selection to D16
$ws.Range("D16").Select()
